$d = $word.ActiveDocument
$sec = $d.Sections(1)

# wdHeaderFooterPrimary = 1 (the "default" header/footer story)
# wdHeaderFooterFirstPage = 2 (the "first page" header/footer story)

function Rename-InlineShape($range, $newName) {
    $paraCount = $range.Paragraphs.Count
    $para = $range.Paragraphs($paraCount)
    $shape = $para.Range.InlineShapes(1)
    $shape.Name = $newName
}

# Headers hold the BTEC logo: image1.jpg -> image2.jpg
Rename-InlineShape $sec.Headers(1).Range "image2.jpg"
Rename-InlineShape $sec.Headers(2).Range "image2.jpg"

# Footers hold the Pearson logo: image2.png -> image1.png
Rename-InlineShape $sec.Footers(1).Range "image1.png"
Rename-InlineShape $sec.Footers(2).Range "image1.png"
